$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the note text for the entry on row 59
$ws.Range("D59").Value = "Finished 1 small problem, worked on a second"

# Update the hours logged for that same entry (0.5 -> 1.25)
$ws.Range("C59").Value = 1.25

$wb.Save()
